$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension/measure metadata cells (row 2)
$ws.Range("B2").Value = "iaest-measure:nucleos-en-el-hogar"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "sdmx-dimension:refArea"

# Update medida/dim swap (row 3)
$ws.Range("B3").Value = "medida"
$ws.Range("D3").Value = "dim"

# Update type/URI cells (row 4)
$ws.Range("B4").Value = "xsd:int"
$ws.Range("D4").Value = "URI-Municipio"
$ws.Range("F4").Value = "URI-Comunidad"

# Remove row 5 (mapping xlsx references no longer needed)
$ws.Rows("5:5").Delete()
